$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 488: new earned SP(1-0-0) value of 1.25 ---
$ws.Range("C488").Value = 1.25

# --- Row 490: add SL(2-0-0) entry, H=2, K gets a textual date note ---
$ws.Range("B490").Value = "SL(2-0-0)"
$ws.Range("H490").Value = 2
$ws.Range("K490").Value = "6/13,14/2023"

# --- Row 491: becomes the "inserted" row: clears the period date, adds SL(1-0-0), H=1, K becomes a real date ---
$ws.Range("B491").Value = "SL(1-0-0)"
$ws.Range("H491").Value = 1
$ws.Range("K488").Copy()
$ws.Range("K491").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K491").Value = 45096

# Shift PERIOD (column A) dates for rows 492:569 down by one row: copy A491:A568 -> A492:A569, then clear A491
$ws.Range("A491:A568").Copy()
$ws.Range("A492:A569").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("A491").ClearContents()

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A8:K597"))

# Row 596 had the footer style/border + formula; duplicate it fully down to the new row 597
# (now inside the resized table), re-adding the structured formula explicitly, then restyle row
# 596 itself as a normal row (copy row 595's look onto it, leaving its own formula/values untouched)
$ws.Range("A596:K596").Copy($ws.Range("A597"))
$ws.Range("G597").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),`"`",Table1[[#This Row],[EARNED]])"
$excel.CutCopyMode = 0

$ws.Range("A595:K595").Copy()
$ws.Range("A596:K596").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$excel.Calculate()
